$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changed from 2023-09-21 (45190) to 2023-09-23 (45192)
# for every existing data row (2..203).
for ($r = 2; $r -le 203; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# Row 203 picks up an explicit row height (15, customHeight) in the edited file.
$ws.Rows.Item(203).RowHeight = 15

# New row 204: a freshly reported avverkningsanmälan.
$ws.Cells.Item(204, 1).Value = "A 44774-2023"
$ws.Cells.Item(204, 2).Value = 45190
$ws.Cells.Item(204, 3).Value = 45192
$ws.Range("B204:C204").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(204, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(204, 5).Value = "BOLLEBYGD"
$ws.Cells.Item(204, 7).Value = 1.6
$ws.Cells.Item(204, 8).Value = 0
$ws.Cells.Item(204, 9).Value = 0
$ws.Cells.Item(204, 10).Value = 0
$ws.Cells.Item(204, 11).Value = 0
$ws.Cells.Item(204, 12).Value = 0
$ws.Cells.Item(204, 13).Value = 0
$ws.Cells.Item(204, 14).Value = 0
$ws.Cells.Item(204, 15).Value = 0
$ws.Cells.Item(204, 16).Value = 0
$ws.Cells.Item(204, 17).Value = 0
$ws.Range("R204").WrapText = $true
